$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "5.0-leche,1.0-huevos,2.0-harinita,1.0-vainilla,"
$ws.Range("C3").Value = "5.0-huevos,1.0-harinita,2.0-manzana,"
$ws.Range("C4").Value = "2.0-huevos,5.0-harinita,1.0-vainilla,"
$ws.Range("C5").Value = "5.0-huevos,5.0-harinita,"
$ws.Range("C6").Value = "1.0-crema,5.0-merengue,2.0-limon,5.0-huevos,4.0-harinita,"
